$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new quotation row (2025-12-08 / serial 45999) below the
# existing data (last row was 93, serial 45998).
$newRow = 94

$ws.Cells.Item($newRow, 1).Value = 45999
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = "15,6471"
$ws.Cells.Item($newRow, 3).Value = "16,0504"
$ws.Cells.Item($newRow, 4).Value = "15,6471"
$ws.Cells.Item($newRow, 5).Value = "15,6471"
